# Darou/Delor yearly balance sheet refresh: database update.
# Each fiscal-year column (D..H) slides one year to the left and a new
# year (1401/12) is appended in column H, matching the new
# "read_price" extraction pass over the source filings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update period headers (row 8) and publish-date headers (row 9)
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("D9").Value = "1399-02-21 (8)"
$ws.Range("E9").Value = "1400-02-29 (9)"
$ws.Range("F9").Value = "1401-03-04 (8)"
$ws.Range("G9").Value = "1402-02-30 (8)"
$ws.Range("H9").Value = "1402-02-30 (2)"

# Update balance sheet data values (shift left one year + append new year column)

# موجودی نقد
$ws.Range("D12").Value = 84699
$ws.Range("E12").Value = 325027
$ws.Range("F12").Value = 609063
$ws.Range("G12").Value = 1316730
$ws.Range("H12").Value = 5072688

# دریافتنی‌های تجاری و سایر دریافتنی‌ها
$ws.Range("D14").Value = 2272069
$ws.Range("E14").Value = 3951373
$ws.Range("F14").Value = 4761789
$ws.Range("G14").Value = 12927260
$ws.Range("H14").Value = 16378616

# موجودی مواد و کالا
$ws.Range("D15").Value = 1397395
$ws.Range("E15").Value = 2183133
$ws.Range("F15").Value = 3726600
$ws.Range("G15").Value = 5452939
$ws.Range("H15").Value = 9271779

# پیش پرداخت ها
$ws.Range("D16").Value = 528453
$ws.Range("E16").Value = 757038
$ws.Range("F16").Value = 1682374
$ws.Range("G16").Value = 1480908
$ws.Range("H16").Value = 3633218

# جمع داراییهای جاری
$ws.Range("D18").Value = 4282616
$ws.Range("E18").Value = 7216571
$ws.Range("F18").Value = 10779826
$ws.Range("G18").Value = 21177837
$ws.Range("H18").Value = 34356301

# حسابها و اسناد دریافتنی تجاری بلند مدت
$ws.Range("D19").Value = 1901
$ws.Range("E19").Value = 1347
$ws.Range("F19").Value = 9045
$ws.Range("G19").Value = 11526
$ws.Range("H19").Value = 24896

# سرمایه گذاریهای بلند مدت (D20 unchanged)
$ws.Range("E20").Value = 24802
$ws.Range("F20").Value = 44502
$ws.Range("G20").Value = 44511
$ws.Range("H20").Value = 93766

# داراییهای ثابت مشهود
$ws.Range("D22").Value = 231786
$ws.Range("E22").Value = 288460
$ws.Range("F22").Value = 558488
$ws.Range("G22").Value = 1127568
$ws.Range("H22").Value = 1657549

# داراییهای نامشهود
$ws.Range("D23").Value = 830
$ws.Range("E23").Value = 1088
$ws.Range("F23").Value = 2458
$ws.Range("G23").Value = 2831
$ws.Range("H23").Value = 1617

# جمع داراییهای غیرجاری
$ws.Range("D26").Value = 245195
$ws.Range("E26").Value = 316348
$ws.Range("F26").Value = 615144
$ws.Range("G26").Value = 1187087
$ws.Range("H26").Value = 1778479

# جمع داراییها
$ws.Range("D27").Value = 4527811
$ws.Range("E27").Value = 7532919
$ws.Range("F27").Value = 11394970
$ws.Range("G27").Value = 22364924
$ws.Range("H27").Value = 36134780

# پرداختنی‌های تجاری و سایر پرداختنی‌ها
$ws.Range("D29").Value = 909288
$ws.Range("E29").Value = 743358
$ws.Range("F29").Value = 1050088
$ws.Range("G29").Value = 2126504
$ws.Range("H29").Value = 3559877

# پیش دریافتها
$ws.Range("D31").Value = 88057
$ws.Range("E31").Value = 108539
$ws.Range("F31").Value = 122131
$ws.Range("G31").Value = 21745
$ws.Range("H31").Value = 26602

# ذخیره مالیات بر درامد
$ws.Range("D32").Value = 185672
$ws.Range("E32").Value = 574542
$ws.Range("F32").Value = 647647
$ws.Range("G32").Value = 1320627
$ws.Range("H32").Value = 1538046

# سود سهام پیشنهادی و پرداختنی
$ws.Range("D33").Value = 78529
$ws.Range("E33").Value = 197983
$ws.Range("F33").Value = 572456
$ws.Range("G33").Value = 1721139
$ws.Range("H33").Value = 4974162

# حصه جاری تسهیلات مالی دریافتی
$ws.Range("D34").Value = 1700646
$ws.Range("E34").Value = 2433521
$ws.Range("F34").Value = 3588014
$ws.Range("G34").Value = 7453965
$ws.Range("H34").Value = 14737294

# جمع بدهیهای جاری
$ws.Range("D37").Value = 2962192
$ws.Range("E37").Value = 4057943
$ws.Range("F37").Value = 5980336
$ws.Range("G37").Value = 12643980
$ws.Range("H37").Value = 24835981

# پیش دریافتهای غیرجاری (D39 becomes the "-" placeholder like the rest of the row)
$ws.Range("D39").Value = "-"

# تسهیلات مالی دریافتی بلند مدت (D40 only)
$ws.Range("D40").Value = 0

# ذخیره مزایای پایان خدمت
$ws.Range("D41").Value = 16586
$ws.Range("E41").Value = 36799
$ws.Range("F41").Value = 62661
$ws.Range("G41").Value = 113925
$ws.Range("H41").Value = 206480

# جمع بدهیهای غیر جاری
$ws.Range("D42").Value = 16586
$ws.Range("E42").Value = 36799
$ws.Range("F42").Value = 62661
$ws.Range("G42").Value = 113925
$ws.Range("H42").Value = 206480

# جمع بدهیهای جاری و غیر جاری
$ws.Range("D43").Value = 2978778
$ws.Range("E43").Value = 4094742
$ws.Range("F43").Value = 6042997
$ws.Range("G43").Value = 12757905
$ws.Range("H43").Value = 25042461

# سرمایه (D45,E45,F45 unchanged)
$ws.Range("G45").Value = 1125000
$ws.Range("H45").Value = 2250000

# وجوه دریافتی بابت افزایش سرمایه (D47,E47,H47 unchanged)
$ws.Range("F47").Value = 653089
$ws.Range("G47").Value = 0

# سهام خزانه (D48,E48 unchanged)
$ws.Range("F48").Value = -57397
$ws.Range("G48").Value = -50268
$ws.Range("H48").Value = -32125

# صرف سهام خزانه (D49 switches from "-" placeholder to numeric 0; E49,F49,G49 unchanged)
$ws.Range("D49").Value = 0
$ws.Range("H49").Value = 86

# اندوخته قانونی (D50,E50,F50 unchanged)
$ws.Range("G50").Value = 112500
$ws.Range("H50").Value = 225000

# مازاد تجدید ارزیابی دارایی های غیر جاری نگه داری شده برای فروش (D52 only)
$ws.Range("D52").Value = "-"

# اندوخته تسعیر ارز داراییها و بدهیهای شرکت های دولتی (D54 only)
$ws.Range("D54").Value = "-"

# سود (زیان) انباشته
$ws.Range("D56").Value = 1054033
$ws.Range("E56").Value = 2943177
$ws.Range("F56").Value = 4261281
$ws.Range("G56").Value = 8419787
$ws.Range("H56").Value = 8649358

# جمع حقوق صاحبان سهام
$ws.Range("D57").Value = 1549033
$ws.Range("E57").Value = 3438177
$ws.Range("F57").Value = 5351973
$ws.Range("G57").Value = 9607019
$ws.Range("H57").Value = 11092319

# جمع بدهیها و حقوق صاحبان سهام
$ws.Range("D58").Value = 4527811
$ws.Range("E58").Value = 7532919
$ws.Range("F58").Value = 11394970
$ws.Range("G58").Value = 22364924
$ws.Range("H58").Value = 36134780
